# Update the "想去人数" (interested count) column F values on both the
# "展览" and "全部类型" worksheets to reflect the newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    2  = 1094
    5  = 4664
    7  = 394
    9  = 922
    11 = 1150
    13 = 636
    15 = 38
    16 = 16
    17 = 276
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
